# "Crear carpeta menús y config separado para datos"
#
# Updates the Tabla1 recipe list on Sheet1:
#   - Smoothie gets an extra ingredient (cacahuete en polvo)
#   - Humus moves from category "legumbres" to "relleno"
#   - Albondigas gets an extra ingredient (tomate)
#   - The placeholder "Prueba2" row becomes the real "Chaofan" recipe
#   - Six brand new recipes are appended (Falafel, Fajitas de pollo,
#     Pan/guarnicion_hidrato, Lentejas/legumbres, Banmian, Ensalada de rucula)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- amend existing rows -----------------------------------------------
$ws.Range("C23").Value = "0, leche; 0, yogur; 0, cacahuete en polvo"

$ws.Range("B27").Value = "relleno"

$ws.Range("C28").Value = "0, carne picada; 2, huevos; 0, harina; 0, pan rallado; 0, tomate"

$ws.Range("A29").Value = "Chaofan"
$ws.Range("C29").Value = "0, arroz; 1, zanahoria; 1, puerro; 1, maíz"

# --- grow the table so the new rows become part of Tabla1 --------------
$lo = $ws.ListObjects.Item("Tabla1")
$lo.Resize($ws.Range("A1:C41"))

# --- append the new recipes ---------------------------------------------
$newRows = @(
    @("Falafel", "legumbres", "1, bolsa de falafel"),
    @("Fajitas de pollo", "carne_blanca", "3, pechuga de pollo; 2, pimiento; 0, cebolla; 0, soja"),
    @("Pan", "guarnicion_hidrato", "0, pan"),
    @("Lentejas", "legumbres", "0, lentejas"),
    @("Banmian", "guarnicion_hidrato", "0, pasta de arroz; 1, puerro; 1, col; 1, carne picada; 0, ajo"),
    @("Ensalada de rúcula", "verdura_cruda", "1, paquete de rúcula; 0, tomates cherry; 0, feta")
)

$row = 33
foreach ($r in $newRows) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row = $row + 1
}

# --- restore the cursor position recorded in the saved view ------------
$ws.Range("B28").Select() | Out-Null
